$wb = $excel.ActiveWorkbook

# sheet1 (sheet index 1)
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(15, 8).Value = 287.34
$ws.Cells.Item(15, 9).Value = 287.34
$ws.Cells.Item(15, 11).Value = 862.02
$ws.Cells.Item(15, 13).Value = -693.02
$ws.Cells.Item(19, 8).Value = 236.45454
$ws.Cells.Item(19, 9).Value = 265.76923
$ws.Cells.Item(19, 10).Value = 194.11111
$ws.Cells.Item(19, 11).Value = 265.76923
$ws.Cells.Item(19, 12).Value = 194.11111
$ws.Cells.Item(19, 13).Value = -90.76922999999999
$ws.Cells.Item(19, 14).Value = -544.1111100000001
$ws.Cells.Item(129, 8).Value = 929.625
$ws.Cells.Item(129, 9).Value = 410
$ws.Cells.Item(129, 10).Value = 1241.4
$ws.Cells.Item(129, 11).Value = 1230
$ws.Cells.Item(129, 12).Value = 3724.2
$ws.Cells.Item(129, 13).Value = 3770
$ws.Cells.Item(129, 14).Value = -13724.2
$ws.Cells.Item(132, 8).Value = 1334.0186
$ws.Cells.Item(132, 9).Value = 1277.6604
$ws.Cells.Item(132, 10).Value = 4321
$ws.Cells.Item(132, 11).Value = 3832.9812
$ws.Cells.Item(132, 12).Value = 12963
$ws.Cells.Item(132, 13).Value = -1302.9812
$ws.Cells.Item(132, 14).Value = -18023
$ws.Cells.Item(137, 8).Value = 984.43243
$ws.Cells.Item(137, 9).Value = 807.1429000000001
$ws.Cells.Item(137, 10).Value = 1536
$ws.Cells.Item(137, 11).Value = 2421.4287
$ws.Cells.Item(137, 12).Value = 4608
$ws.Cells.Item(137, 13).Value = 128.5712999999996
$ws.Cells.Item(137, 14).Value = -9708

# sheet2 (sheet index 2)
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(2, 8).Value = 2454.1333
$ws.Cells.Item(2, 9).Value = 1004
$ws.Cells.Item(2, 10).Value = 4111.4287
$ws.Cells.Item(2, 11).Value = 1004
$ws.Cells.Item(2, 12).Value = 4111.4287
$ws.Cells.Item(2, 13).Value = -891
$ws.Cells.Item(2, 14).Value = -4337.4287
$ws.Cells.Item(37, 8).Value = 8200.75
$ws.Cells.Item(37, 10).Value = 8200.75
$ws.Cells.Item(37, 12).Value = 8200.75
$ws.Cells.Item(37, 14).Value = -8746.75
$ws.Cells.Item(45, 8).Value = 2825.257
$ws.Cells.Item(45, 9).Value = 2681.9565
$ws.Cells.Item(45, 11).Value = 2681.9565
$ws.Cells.Item(45, 13).Value = -2304.9565
$ws.Cells.Item(80, 8).Value = 16195.714
$ws.Cells.Item(80, 10).Value = 18495
$ws.Cells.Item(80, 12).Value = 18495
$ws.Cells.Item(80, 14).Value = -20491
$ws.Cells.Item(83, 8).Value = 16195.714
$ws.Cells.Item(83, 10).Value = 18495
$ws.Cells.Item(83, 12).Value = 55485
$ws.Cells.Item(83, 14).Value = -65469
$ws.Cells.Item(97, 8).Value = 1816.0769
$ws.Cells.Item(97, 9).Value = 1805.7142
$ws.Cells.Item(97, 10).Value = 1828.1666
$ws.Cells.Item(97, 11).Value = 1805.7142
$ws.Cells.Item(97, 12).Value = 1828.1666
$ws.Cells.Item(97, 13).Value = -1309.7142
$ws.Cells.Item(97, 14).Value = -2820.1666
$ws.Cells.Item(102, 8).Value = 5120
$ws.Cells.Item(102, 9).Value = 3980
$ws.Cells.Item(102, 10).Value = 7400
$ws.Cells.Item(102, 11).Value = 3980
$ws.Cells.Item(102, 12).Value = 7400
$ws.Cells.Item(102, 13).Value = -2358
$ws.Cells.Item(102, 14).Value = -10644
$ws.Cells.Item(116, 8).Value = 2454.1333
$ws.Cells.Item(116, 9).Value = 1004
$ws.Cells.Item(116, 10).Value = 4111.4287
$ws.Cells.Item(116, 11).Value = 1004
$ws.Cells.Item(116, 12).Value = 4111.4287
$ws.Cells.Item(116, 13).Value = 1290
$ws.Cells.Item(116, 14).Value = -8699.4287
$ws.Cells.Item(122, 8).Value = 2099.5945
$ws.Cells.Item(122, 9).Value = 2082.2942
$ws.Cells.Item(122, 10).Value = 2295.6667
$ws.Cells.Item(122, 11).Value = 6246.882599999999
$ws.Cells.Item(122, 12).Value = 6887.000100000001
$ws.Cells.Item(122, 13).Value = -3796.882599999999
$ws.Cells.Item(122, 14).Value = -11787.0001
$ws.Cells.Item(125, 8).Value = 50000
$ws.Cells.Item(125, 10).Value = 50000
$ws.Cells.Item(125, 12).Value = 50000
$ws.Cells.Item(125, 14).Value = -59840

# sheet3 (sheet index 3)
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(3, 8).Value = 2454.1333
$ws.Cells.Item(3, 9).Value = 1004
$ws.Cells.Item(3, 10).Value = 4111.4287
$ws.Cells.Item(3, 11).Value = 1004
$ws.Cells.Item(3, 12).Value = 4111.4287
$ws.Cells.Item(3, 13).Value = -890
$ws.Cells.Item(3, 14).Value = -4339.4287
$ws.Cells.Item(86, 8).Value = 2419.9546
$ws.Cells.Item(86, 9).Value = 2315.9143
$ws.Cells.Item(86, 10).Value = 2537.4194
$ws.Cells.Item(86, 11).Value = 2315.9143
$ws.Cells.Item(86, 12).Value = 2537.4194
$ws.Cells.Item(86, 13).Value = -1192.9143
$ws.Cells.Item(86, 14).Value = -4783.419400000001
$ws.Cells.Item(89, 8).Value = 2419.9546
$ws.Cells.Item(89, 9).Value = 2315.9143
$ws.Cells.Item(89, 10).Value = 2537.4194
$ws.Cells.Item(89, 11).Value = 11579.5715
$ws.Cells.Item(89, 12).Value = 12687.097
$ws.Cells.Item(89, 13).Value = -5963.5715
$ws.Cells.Item(89, 14).Value = -23919.097
$ws.Cells.Item(99, 8).Value = 1938.0526
$ws.Cells.Item(99, 9).Value = 1507.1818
$ws.Cells.Item(99, 10).Value = 2530.5
$ws.Cells.Item(99, 11).Value = 1507.1818
$ws.Cells.Item(99, 12).Value = 2530.5
$ws.Cells.Item(99, 13).Value = -9.181800000000067
$ws.Cells.Item(99, 14).Value = -5526.5
$ws.Cells.Item(105, 8).Value = 1571.2826
$ws.Cells.Item(105, 9).Value = 1488.5358
$ws.Cells.Item(105, 11).Value = 1488.5358
$ws.Cells.Item(105, 13).Value = 258.4641999999999
$ws.Cells.Item(107, 8).Value = 770466.3
$ws.Cells.Item(107, 9).Value = 1188.3
$ws.Cells.Item(107, 10).Value = 3334726.2
$ws.Cells.Item(107, 11).Value = 1188.3
$ws.Cells.Item(107, 12).Value = 3334726.2
$ws.Cells.Item(107, 13).Value = 731.7
$ws.Cells.Item(107, 14).Value = -3338566.2
$ws.Cells.Item(124, 8).Value = 27000
$ws.Cells.Item(124, 10).Value = 27000
$ws.Cells.Item(124, 12).Value = 27000
$ws.Cells.Item(124, 14).Value = -36820

# sheet4 (sheet index 4)
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(50, 8).Value = 8320.799999999999
$ws.Cells.Item(50, 10).Value = 8320.799999999999
$ws.Cells.Item(50, 12).Value = 8320.799999999999
$ws.Cells.Item(50, 14).Value = -9570.799999999999
$ws.Cells.Item(59, 8).Value = 9411.4
$ws.Cells.Item(59, 9).Value = 4552
$ws.Cells.Item(59, 10).Value = 12651
$ws.Cells.Item(59, 11).Value = 4552
$ws.Cells.Item(59, 12).Value = 12651
$ws.Cells.Item(59, 13).Value = -3407
$ws.Cells.Item(59, 14).Value = -14941
$ws.Cells.Item(74, 8).Value = 17925.5
$ws.Cells.Item(74, 10).Value = 17925.5
$ws.Cells.Item(74, 12).Value = 17925.5
$ws.Cells.Item(74, 14).Value = -19673.5
$ws.Cells.Item(77, 8).Value = 17925.5
$ws.Cells.Item(77, 10).Value = 17925.5
$ws.Cells.Item(77, 12).Value = 53776.5
$ws.Cells.Item(77, 14).Value = -62512.5
$ws.Cells.Item(80, 8).Value = 23450
$ws.Cells.Item(80, 10).Value = 23450
$ws.Cells.Item(80, 12).Value = 23450
$ws.Cells.Item(80, 14).Value = -25696
$ws.Cells.Item(83, 8).Value = 23450
$ws.Cells.Item(83, 10).Value = 23450
$ws.Cells.Item(83, 12).Value = 70350
$ws.Cells.Item(83, 14).Value = -81582
$ws.Cells.Item(111, 8).Value = 30000
$ws.Cells.Item(111, 10).Value = 30000
$ws.Cells.Item(111, 12).Value = 30000
$ws.Cells.Item(111, 14).Value = -38180

# sheet5 (sheet index 5)
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(5, 8).Value = 739.3
$ws.Cells.Item(5, 9).Value = 624.1875
$ws.Cells.Item(5, 10).Value = 1199.75
$ws.Cells.Item(5, 11).Value = 1872.5625
$ws.Cells.Item(5, 12).Value = 3599.25
$ws.Cells.Item(5, 13).Value = -1760.5625
$ws.Cells.Item(5, 14).Value = -3823.25
$ws.Cells.Item(131, 8).Value = 5051398.5
$ws.Cells.Item(131, 9).Value = 1222.4166
$ws.Cells.Item(131, 10).Value = 5747975
$ws.Cells.Item(131, 11).Value = 3667.2498
$ws.Cells.Item(131, 12).Value = 17243925
$ws.Cells.Item(131, 13).Value = 1372.7502
$ws.Cells.Item(131, 14).Value = -17254005
$ws.Cells.Item(135, 8).Value = 739.3
$ws.Cells.Item(135, 9).Value = 624.1875
$ws.Cells.Item(135, 10).Value = 1199.75
$ws.Cells.Item(135, 11).Value = 5617.6875
$ws.Cells.Item(135, 12).Value = 10797.75
$ws.Cells.Item(135, 13).Value = -3082.6875
$ws.Cells.Item(135, 14).Value = -15867.75

# sheet6 (sheet index 6)
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(97, 8).Value = 1547.7778
$ws.Cells.Item(97, 9).Value = 883.6667
$ws.Cells.Item(97, 10).Value = 2876
$ws.Cells.Item(97, 11).Value = 883.6667
$ws.Cells.Item(97, 12).Value = 2876
$ws.Cells.Item(97, 13).Value = -387.6667
$ws.Cells.Item(97, 14).Value = -3868
$ws.Cells.Item(102, 8).Value = 1505.5927
$ws.Cells.Item(102, 9).Value = 1304.05
$ws.Cells.Item(102, 10).Value = 2081.4285
$ws.Cells.Item(102, 11).Value = 1304.05
$ws.Cells.Item(102, 12).Value = 2081.4285
$ws.Cells.Item(102, 13).Value = 317.95
$ws.Cells.Item(102, 14).Value = -5325.4285
$ws.Cells.Item(113, 8).Value = 55557200
$ws.Cells.Item(113, 9).Value = 1680.1428
$ws.Cells.Item(113, 10).Value = 90910710
$ws.Cells.Item(113, 11).Value = 1680.1428
$ws.Cells.Item(113, 12).Value = 90910710
$ws.Cells.Item(113, 13).Value = 489.8571999999999
$ws.Cells.Item(113, 14).Value = -90915050
$ws.Cells.Item(123, 8).Value = 13103.25
$ws.Cells.Item(123, 10).Value = 13103.25
$ws.Cells.Item(123, 12).Value = 13103.25
$ws.Cells.Item(123, 14).Value = -18003.25
$ws.Cells.Item(126, 8).Value = 12822473
$ws.Cells.Item(126, 9).Value = 2549.6667
$ws.Cells.Item(126, 10).Value = 23810978
$ws.Cells.Item(126, 11).Value = 7649.000100000001
$ws.Cells.Item(126, 12).Value = 71432934
$ws.Cells.Item(126, 13).Value = -5179.000100000001
$ws.Cells.Item(126, 14).Value = -71437874

# sheet7 (sheet index 7)
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(36, 8).Value = 33238.332
$ws.Cells.Item(36, 10).Value = 33238.332
$ws.Cells.Item(36, 12).Value = 33238.332
$ws.Cells.Item(36, 14).Value = -34362.332
